$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the `_GoBack` bookmark that used to sit after "...management
#    system." (it is being relocated to the end of the new "Apache 2" line
#    further down in the skills list).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2. Locate the "Backend services" skills block so we can insert the new
#    "Nginx" and "Apache 2" lines around the existing entries.
# ---------------------------------------------------------------------------
$headingIndex = 0
$oauthIndex = 0
$expressIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Backend services`r") {
        $headingIndex = $i
    }
    if ($headingIndex -gt 0 -and $oauthIndex -eq 0 -and $t -match "OAuth \(custom flow") {
        $oauthIndex = $i
    }
    if ($oauthIndex -gt 0 -and $expressIndex -eq 0 -and $t -match "ExpressJS") {
        $expressIndex = $i
    }
    if ($expressIndex -gt 0) {
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Insert the new "Nginx" paragraph right after the "Backend services"
#    heading, i.e. right before the existing OAuth paragraph.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingPara.Range.InsertParagraphAfter() | Out-Null

$nginxPara = $d.Paragraphs.Item($headingIndex + 1)
$nginxXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="964" w:hanging="964"/><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>&#9733;&#9733;&#9733;&#9733;&#9734;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Nginx (</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">proxy, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>ws</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>fastcgi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>tsl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>ssl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>, templating, dynamic)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$nginxPara.Range.InsertXML($nginxXml)

# ---------------------------------------------------------------------------
# 4. Insert the new "Apache 2" paragraph right after the existing "ExpressJS"
#    paragraph (which shifted down by one because of the insertion above).
# ---------------------------------------------------------------------------
$expressPara = $d.Paragraphs.Item($expressIndex + 1)
$expressPara.Range.InsertParagraphAfter() | Out-Null

$apachePara = $d.Paragraphs.Item($expressIndex + 2)
$apacheXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="964" w:hanging="964"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t>&#9733;&#9733;&#9734;&#9734;&#9734;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Apache 2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$apachePara.Range.InsertXML($apacheXml)

# Re-add the `_GoBack` bookmark as a zero-width bookmark at the very end of
# the new "Apache 2" paragraph's text (right after "Apache 2", before the
# paragraph mark).
$apacheParaAfter = $d.Paragraphs.Item($expressIndex + 2)
$bmRange = $apacheParaAfter.Range.Duplicate
$bmRange.SetRange($apacheParaAfter.Range.End - 1, $apacheParaAfter.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 5. Add a `lastRenderedPageBreak` marker to the final bullet point ("Active
#    member and elected Vice-chairman...").
# ---------------------------------------------------------------------------
$lastIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Active member and elected Vice-chairman") {
        $lastIndex = $i
    }
}
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t>Active member and elected Vice-chairman (2016-2017) of Academic Catholic Student Association Soli Deo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$lastPara.Range.InsertXML($lastXml)

Write-Host "Done."
